$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.806.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.767.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.26%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.561"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.50%  "
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.201.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.790.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.929"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.741.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0974"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "273.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.145"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.34%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "51.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0462"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.78%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("E40").Value = "  -4.83%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.114"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.51"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("E43").Value = "  -2.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "121.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.059.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.923"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.43%  "
$ws.Range("E51").Value = "  -0.27%  "
